$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '43.080.59'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.45%  '

# Row 3: Ethereum -> Ethereum
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '2.377.05'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +2.21%  '

# Row 4: TetherUSD -> TetherUSD
$ws.Range("E4").Value = '  +0.10%  '

# Row 5: BNB -> BNB
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '302.04'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.20%  '

# Row 6: Solana -> Solana
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '96.48'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -0.05%  '

# Row 7: XRP -> XRP
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.504'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.66%  '

# Row 8: USDC -> USDC
$ws.Range("E8").Value = '  -0.08%  '

# Row 9: Cardano -> Cardano
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.499'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.36%  '

# Row 10: Avalanche -> Avalanche
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '34.25'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -1.22%  '

# Row 11: Dogecoin -> Dogecoin
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.0788'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.31%  '

# Row 12: TRON -> TRON
$ws.Range("E12").Value = '  +2.36%  '

# Row 13: Chainlink -> Chainlink
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '18.22'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -4.96%  '

# Row 14: Polkadot -> Polkadot
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '6.80'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +0.31%  '

# Row 15: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '2.747.92'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +2.15%  '

# Row 16: WrappedEther -> WrappedEther
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '2.437.19'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +4.75%  '

# Row 17: Polygon -> Polygon
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.804'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +1.65%  '

# Row 18: WrappedBTC -> WrappedBTC
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '43.045.86'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.53%  '

# Row 19: InternetComputer(DFINITY) -> InternetComputer(DFINITY)
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.16'
$c.Style = "Normal"
$ws.Range("E19").Value = '  -1.33%  '

# Row 20: Uniswap -> Uniswap
$ws.Range("E20").Value = '  +2.57%  '

# Row 21: ShibaInu -> ShibaInu
$ws.Range("E21").Value = '  -0.52%  '

# Row 22: Litecoin -> Litecoin
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '68.15'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '

# Row 23: BitcoinCash -> BitcoinCash
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '235.36'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.60%  '

# Row 24: ImmutableX -> ImmutableX
$ws.Range("E24").Value = '  -1.54%  '

# Row 25: Dai -> PancakeSwap
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.44'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.44%  '

# Row 26: PancakeSwap -> Dai
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.06%  '

# Row 27: EthereumClassic -> EthereumClassic
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '24.92'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.77%  '

# Row 28: Toncoin -> Toncoin
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.37'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -0.25%  '

# Row 29: Cosmos -> Cosmos
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.25'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +1.16%  '

# Row 30: InjectiveProtocol -> InjectiveProtocol
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '31.49'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -2.90%  '

# Row 31: FirstDigitalUSD -> FirstDigitalUSD
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.04%  '

# Row 32: Filecoin -> Filecoin
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '5.07'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.84%  '

# Row 33: Hedera -> Hedera
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.0739'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +5.25%  '

# Row 34: Celestia -> Celestia
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '17.55'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -1.31%  '

# Row 35: Kaspa -> Kaspa
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '0.105'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +5.50%  '

# Row 36: ARBITRUM -> ARBITRUM
$ws.Range("E36").Value = '  +4.45%  '

# Row 37: RenderToken -> RenderToken
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '4.36'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -2.70%  '

# Row 38: WEMIXToken -> WEMIXToken
$ws.Range("E38").Value = '  -0.92%  '

# Row 39: LidoDAOToken -> LidoDAOToken
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.80'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +2.14%  '

# Row 40: EnergySwap -> EnergySwap
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '22.29'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +8.53%  '

# Row 41: Stellar -> Stellar
$ws.Range("E41").Value = '  -0.68%  '

# Row 42: Monero -> Monero
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '116.06'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -30.32%  '

# Row 43: Maker -> Maker
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.953.13'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.50%  '

# Row 44: VeChain -> VeChain
$ws.Range("E44").Value = '  +0.39%  '

# Row 45: ApeXProtocol -> ApeXProtocol
$ws.Range("E45").Value = '  +2.06%  '

# Row 46: NEARProtocol -> NEARProtocol
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.74'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.79%  '

# Row 47: FraxShare -> FraxShare
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '9.15'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -11.51%  '

# Row 48: Stacks -> RocketPoolETH
$ws.Range("B48").Value = 'RocketPoolETH'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.605.24'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +1.86%  '

# Row 49: MultiversX -> Stacks
$ws.Range("B49").Value = 'Stacks'
$ws.Range("C49").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.52'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +2.45%  '

# Row 50: BitcoinSV -> MultiversX
$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '52.37'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.16%  '

# Row 51: TrustWalletToken -> BitcoinSV
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '72.21'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.04%  '
